# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh the
# handoff timestamps on the Overview / zh-cn / de-de sheets, then resize the
# Status columns to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-31 00:41:58"

# --- zh-cn sheet --------------------------------------------------------
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("H2").Value = "2016-08-31 00:41:53"

# --- de-de sheet --------------------------------------------------------
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("H2").Value = "2016-08-31 00:41:58"

# --- Resize the Status columns to fit the new text -----------------------
# (target authored width is 17.2159881591797; the interop layer quantizes
# ColumnWidth writes to 1/6-character steps, so feed it the input that lands
# on the nearest attainable grid point.)
$ws1.Range("E1:F1").ColumnWidth = 16.333333333333332
$ws2.Range("C1").ColumnWidth = 16.333333333333332
$ws3.Range("C1").ColumnWidth = 16.333333333333332
